$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as plain text so that numeric-looking
    # strings (e.g. "0.500", "19.60") are not silently converted into
    # actual numbers (which would drop formatting such as trailing zeros).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "27.114.59"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.647.27"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.54%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "219.97"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.500"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.61%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.254"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.0626"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "19.60"
$ws.Range("E10").Value = "  +2.05%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0848"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.875.45"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13 - now WrappedEther (was Polkadot)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "1.649.52"
$ws.Range("E13").Value = "  +0.69%  "

# Row 14 - now Polkadot (was WrappedEther)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "4.20"
$ws.Range("E14").Value = "  +0.87%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.532"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "66.21"
$ws.Range("E16").Value = "  +2.37%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "27.042.87"
$ws.Range("E17").Value = "  +1.00%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.14%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "220.56"
$ws.Range("E19").Value = "  +2.91%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.39%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "4.42"
$ws.Range("E21").Value = "  +0.61%  "

# Row 22 - Chainlink
Set-TextValue $ws.Range("D22") "6.65"
$ws.Range("E22").Value = "  +6.41%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -2.20%  "

# Row 24 - Avalanche
Set-TextValue $ws.Range("D24") "9.26"
$ws.Range("E24").Value = "  -1.11%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "147.72"
$ws.Range("E25").Value = "  +0.61%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.49%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "7.43"
$ws.Range("E27").Value = "  +3.67%  "

# Row 28 - Stellar
Set-TextValue $ws.Range("D28") "0.119"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "15.87"
$ws.Range("E29").Value = "  +1.42%  "

# Row 30 - Hedera
Set-TextValue $ws.Range("D30") "0.0514"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.33%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.60%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.20%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +2.27%  "

# Row 35 - Maker
Set-TextValue $ws.Range("D35") "1.265.63"
$ws.Range("E35").Value = "  -1.96%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.16%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.0174"
$ws.Range("E37").Value = "  -1.68%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -0.54%  "

# Row 39 - ARBITRUM
Set-TextValue $ws.Range("D39") "0.827"
$ws.Range("E39").Value = "  +0.57%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.36%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  +1.15%  "

# Row 43 - RocketPoolETH
Set-TextValue $ws.Range("D43") "1.786.46"
$ws.Range("E43").Value = "  +0.17%  "

# Row 44 - Aave
Set-TextValue $ws.Range("D44") "61.97"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45 - Quant
Set-TextValue $ws.Range("D45") "92.67"
$ws.Range("E45").Value = "  +1.21%  "

# Row 46 - MXToken
$ws.Range("E46").Value = "  -7.80%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -0.24%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  -0.85%  "

# Row 49 - EnergySwap
Set-TextValue $ws.Range("D49") "7.65"
$ws.Range("E49").Value = "  +0.01%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +0.12%  "

# Row 51 - now Mantle (was USDD)
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.406"
$ws.Range("E51").Value = "  -0.15%  "
